$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.503.74'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '3.812.60'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'613.34"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('D6').Value = "'163.43"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.61%  '
$ws.Range('D7').Value = '3.812.65'
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = "'0.517"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.33%  '
$ws.Range('D10').Value = "'0.160"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('D11').Value = "'0.449"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = "'6.81"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.52%  '
$ws.Range('E13').Value = '  -1.92%  '
$ws.Range('D14').Value = "'35.06"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.88%  '
$ws.Range('D15').Value = '4.442.94'
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('D16').Value = '3.788.86'
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('D17').Value = '68.449.66'
$ws.Range('E17').Value = '  +0.89%  '
$ws.Range('D18').Value = "'18.07"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.76%  '
$ws.Range('B19').Value = 'Polkadot'
$ws.Range('C19').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Value = "'7.06"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('B20').Value = 'TRON'
$ws.Range('C20').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D20').Value = "'0.113"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = "'462.79"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('D22').Value = "'9.61"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.86%  '
$ws.Range('D23').Value = "'0.697"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.77%  '
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('D25').Value = "'83.34"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.01%  '
$ws.Range('D26').Value = "'11.95"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.39%  '
$ws.Range('D27').Value = "'2.10"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').Value = "'9.96"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').Value = '3.951.33'
$ws.Range('E30').Value = '  -0.26%  '
$ws.Range('D31').Value = "'2.62"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.30%  '
$ws.Range('D32').Value = "'2.20"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('D33').Value = "'7.20"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.80%  '
$ws.Range('D34').Value = "'28.88"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.46%  '
$ws.Range('D35').Value = "'1.00"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = "'9.03"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.63%  '
$ws.Range('D37').Value = "'0.100"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('E38').Value = '  +5.19%  '
$ws.Range('D39').Value = "'5.87"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.88%  '
$ws.Range('D40').Value = "'0.982"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.79%  '
$ws.Range('D41').Value = "'0.998"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').Value = "'3.09"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('D44').Value = "'153.19"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.02%  '
$ws.Range('D45').Value = "'42.90"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.66%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').Value = "'0.295"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.87%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D47').Value = "'1.40"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.93%  '
$ws.Range('D48').Value = "'46.61"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('D49').Value = "'8.35"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').Value = "'376.99"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.36%  '
